$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.39000491208574
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("E3").Value = "Colleen"
$ws.Range("G3").Value = 6.091303748649244
$ws.Range("H3").Value = "White"
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("E4").Value = "Jewel"
$ws.Range("G4").Value = 6.07244745832581
$ws.Range("H4").Value = "Black or African American"
$ws.Range("C5").Value = 34
$ws.Range("D5").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("E5").Value = "Tina"
$ws.Range("G5").Value = 5.476255900907384
$ws.Range("H5").Value = "White"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E6").Value = "Annes"
$ws.Range("G6").Value = 5.199586145078674
$ws.Range("H6").Value = "Asian"
$ws.Range("G7").Value = 4.182434273974512
$ws.Range("G8").Value = 1.405244066240008
$ws.Range("G9").Value = 1.285989741820513
$ws.Range("G10").Value = 0.4301310047900727
$ws.Range("G11").Value = 0.3677338533072753
$ws.Range("C12").Value = 33
$ws.Range("D12").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("E12").Value = "Shaniek"
$ws.Range("G12").Value = 0.2397686769137523
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("E13").Value = "Shadaisia"
$ws.Range("G13").Value = 0.2362855029629706
$ws.Range("G14").Value = 13.08405170480421
$ws.Range("G15").Value = 8.317598354541474
$ws.Range("G16").Value = 7.382070019746715
$ws.Range("G17").Value = 7.27978505289396
$ws.Range("C18").Value = 32
$ws.Range("D18").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E18").Value = "Jamarii"
$ws.Range("G18").Value = 5.430038462157364
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("E19").Value = "Juan"
$ws.Range("G19").Value = 5.411470426993446
$ws.Range("H19").Value = "Hispanic"
$ws.Range("C20").Value = 22
$ws.Range("D20").Value = "60db4fde6193c50664c9c478"
$ws.Range("E20").Value = "Edosagbe"
$ws.Range("G20").Value = 5.05463621270477
$ws.Range("H20").Value = "Black or African American"
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = "5e2522d6b734b47915f88275"
$ws.Range("E21").Value = "Corey"
$ws.Range("G21").Value = 4.430868679986358
$ws.Range("C22").Value = 33
$ws.Range("D22").Value = "60b322994d0b901954690036"
$ws.Range("E22").Value = "Brennan"
$ws.Range("G22").Value = 4.412626648038093
$ws.Range("G23").Value = 3.327095999247362
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = "6097b95056caf5ebb2720002"
$ws.Range("E24").Value = "Damian"
$ws.Range("G24").Value = 2.26493513038394
$ws.Range("H24").Value = "Black or African American"
$ws.Range("C25").Value = 29
$ws.Range("D25").Value = "60b83826821417f8e484a207"
$ws.Range("E25").Value = "Eli"
$ws.Range("G25").Value = 2.202334476874346
$ws.Range("H25").Value = "White"
